$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 200
$ws.Range("I39").Value = 200
$ws.Range("K39").Value = 600
$ws.Range("M39").Value = -304

$ws.Range("H100").Value = 2125.7827
$ws.Range("I100").Value = 1840.4166
$ws.Range("J100").Value = 2437.0908
$ws.Range("K100").Value = 1840.4166
$ws.Range("L100").Value = 2437.0908
$ws.Range("M100").Value = -1299.4166
$ws.Range("N100").Value = -3519.0908

$ws.Range("H125").Value = 2028.8889
$ws.Range("I125").Value = 2032
$ws.Range("J125").Value = 2028
$ws.Range("K125").Value = 18288
$ws.Range("L125").Value = 18252
$ws.Range("M125").Value = -15828
$ws.Range("N125").Value = -23172

$ws.Range("H129").Value = 1004.4889
$ws.Range("I129").Value = 1543.3
$ws.Range("J129").Value = 937.1375
$ws.Range("K129").Value = 4629.9
$ws.Range("L129").Value = 2811.4125
$ws.Range("M129").Value = 370.1000000000004
$ws.Range("N129").Value = -12811.4125

$ws.Range("H135").Value = 93752640
$ws.Range("I135").Value = 90910530
$ws.Range("J135").Value = 100005300
$ws.Range("K135").Value = 818194770
$ws.Range("L135").Value = 900047700
$ws.Range("M135").Value = -818192235
$ws.Range("N135").Value = -900052770

$ws.Range("H138").Value = 1690.2526
$ws.Range("I138").Value = 1243.4117
$ws.Range("J138").Value = 1782.8903
$ws.Range("K138").Value = 3730.2351
$ws.Range("L138").Value = 5348.6709
$ws.Range("M138").Value = 1409.7649
$ws.Range("N138").Value = -15628.6709


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2046.174
$ws.Range("I2").Value = 2002.8182
$ws.Range("K2").Value = 2002.8182
$ws.Range("M2").Value = -1889.8182

$ws.Range("H116").Value = 2046.174
$ws.Range("I116").Value = 2002.8182
$ws.Range("K116").Value = 2002.8182
$ws.Range("M116").Value = 291.1818000000001

$ws.Range("H122").Value = 2344.5264
$ws.Range("I122").Value = 2201.6365
$ws.Range("J122").Value = 2541
$ws.Range("K122").Value = 6604.9095
$ws.Range("L122").Value = 7623
$ws.Range("M122").Value = -4154.9095
$ws.Range("N122").Value = -12523


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2046.174
$ws.Range("I3").Value = 2002.8182
$ws.Range("K3").Value = 2002.8182
$ws.Range("M3").Value = -1888.8182

$ws.Range("H8").Value = 5633.3335
$ws.Range("J8").Value = 7950
$ws.Range("L8").Value = 7950
$ws.Range("N8").Value = -8230

$ws.Range("H105").Value = 3022.3044
$ws.Range("I105").Value = 1638.2
$ws.Range("J105").Value = 4087
$ws.Range("K105").Value = 1638.2
$ws.Range("L105").Value = 4087
$ws.Range("M105").Value = 108.8
$ws.Range("N105").Value = -7581

$ws.Range("H134").Value = 2441.4
$ws.Range("I134").Value = 1508.4857
$ws.Range("J134").Value = 3374.3142
$ws.Range("K134").Value = 4525.4571
$ws.Range("L134").Value = 10122.9426
$ws.Range("M134").Value = -1990.4571
$ws.Range("N134").Value = -15192.9426


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3140.1304
$ws.Range("I99").Value = 3264.6
$ws.Range("J99").Value = 3105.5557
$ws.Range("K99").Value = 3264.6
$ws.Range("L99").Value = 3105.5557
$ws.Range("M99").Value = -1766.6
$ws.Range("N99").Value = -6101.5557

$ws.Range("H107").Value = 542.6818
$ws.Range("I107").Value = 393.7857
$ws.Range("J107").Value = 803.25
$ws.Range("K107").Value = 393.7857
$ws.Range("L107").Value = 803.25
$ws.Range("M107").Value = 1526.2143
$ws.Range("N107").Value = -4643.25

$ws.Range("H122").Value = 87580.57000000001
$ws.Range("I122").Value = 121951.4
$ws.Range("J122").Value = 1653.5
$ws.Range("K122").Value = 365854.2
$ws.Range("L122").Value = 4960.5
$ws.Range("M122").Value = -363404.2
$ws.Range("N122").Value = -9860.5

$ws.Range("H126").Value = 3140.1304
$ws.Range("I126").Value = 3264.6
$ws.Range("J126").Value = 3105.5557
$ws.Range("K126").Value = 9793.799999999999
$ws.Range("L126").Value = 9316.667099999999
$ws.Range("M126").Value = -7323.799999999999
$ws.Range("N126").Value = -14256.6671

$ws.Range("H134").Value = 739546.25
$ws.Range("I134").Value = 1007529.4
$ws.Range("K134").Value = 3022588.2
$ws.Range("M134").Value = -3020053.2


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2000
$ws.Range("J48").Value = 2000
$ws.Range("L48").Value = 6000
$ws.Range("N48").Value = -6500


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1635.125
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

$ws.Range("H122").Value = 1742.8572
$ws.Range("I122").Value = 1742.8572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5228.571599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2778.571599999999
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 23812302
$ws.Range("I132").Value = 52633410
$ws.Range("J132").Value = 3562.1738
$ws.Range("K132").Value = 157900230
$ws.Range("L132").Value = 10686.5214
$ws.Range("M132").Value = -157897700
$ws.Range("N132").Value = -15746.5214


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 19502.5
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 29005
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 29005
$ws.Range("M14").Value = -9828
$ws.Range("N14").Value = -29349

$ws.Range("H40").Value = 3750.25
$ws.Range("I40").Value = 3750.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3750.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3614.25
$ws.Range("N40").ClearContents()

$ws.Range("I122").Value = 102400.4
$ws.Range("K122").Value = 307201.2
$ws.Range("M122").Value = -304751.2

$ws.Range("H132").Value = 2868.4
$ws.Range("I132").Value = 2350.182
$ws.Range("J132").Value = 3745.3845
$ws.Range("K132").Value = 7050.545999999999
$ws.Range("L132").Value = 11236.1535
$ws.Range("M132").Value = -4520.545999999999
$ws.Range("N132").Value = -16296.1535


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1176826.8
$ws.Range("I132").Value = 1891643.6
$ws.Range("J132").Value = 2484.7856
$ws.Range("K132").Value = 5674930.800000001
$ws.Range("L132").Value = 7454.3568
$ws.Range("M132").Value = -5672400.800000001
$ws.Range("N132").Value = -12514.3568

